# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation" on all
#   sheets that surface it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Shrink the now-narrower "status" columns to match the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status values (same shared text on every sheet).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the status columns for the shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
